$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: FAC 890, Hospital Woodstock, Name (blank), Title (blank), Pattern "Manual entry D"
$ws.Range("A4").Value = 890
$ws.Range("B4").Value = "Woodstock"
$ws.Range("E4").Value = "Manual entry D"

# Row 5: FAC 896, Hospital Red Lake, Name Angela Bishop, Title Interim President, Pattern "h2_name and p title"
$ws.Range("A5").Value = 896
$ws.Range("B5").Value = "Red Lake"
$ws.Range("C5").Value = "Angela Bishop"
$ws.Range("D5").Value = "Interim President"
$ws.Range("E5").Value = "h2_name and p title"

# Row 6: FAC 942, Hospital Hamilton Heatl, Name Rochelle Reid, Title Stratetic lead, Pattern "div_classes"
$ws.Range("A6").Value = 942
$ws.Range("B6").Value = "Hamilton Heatl"
$ws.Range("C6").Value = "Rochelle Reid"
$ws.Range("D6").Value = "Stratetic lead"
$ws.Range("E6").Value = "div_classes"

# Update the active selection to match the edited workbook state
$ws.Range("C7").Select()
